# "la til skoleeksamen V22" - add the Spring 2022 ("2022 - Vår") school exam
# row to the exam-archive overview sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row right after the existing last row (row 28 -> row 29).
$ws.Range("A29").Value = "2022 - Vår"
$ws.Range("B29").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-22-v.pdf)"
$ws.Range("C29").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-22-v-fasit.pdf)"

# Match the author's final selection/view after adding the row.
[void]$ws.Range("C30").Select()
